$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to Price (D) column cells being updated so that
# numeric-looking strings (e.g. "7.14", "0.870") are preserved exactly as text
# instead of being auto-converted to numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '41.188.62'
$ws.Range('E2').Value = '  +1.26%  '

$ws.Range('D3').Value = '2.181.93'
$ws.Range('E3').Value = '  +0.19%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '254.33'
$ws.Range('E5').Value = '  +6.98%  '

$ws.Range('D6').Value = '0.627'
$ws.Range('E6').Value = '  +1.79%  '

$ws.Range('D7').Value = '67.85'
$ws.Range('E7').Value = '  -2.37%  '

$ws.Range('E8').Value = '  +0.10%  '

$ws.Range('D9').Value = '0.579'
$ws.Range('E9').Value = '  +9.01%  '

$ws.Range('D10').Value = '37.55'
$ws.Range('E10').Value = '  +3.85%  '

$ws.Range('D11').Value = '58.72'
$ws.Range('E11').Value = '  +2.16%  '

$ws.Range('D12').Value = '0.0934'
$ws.Range('E12').Value = '  +0.15%  '

$ws.Range('D13').Value = '7.14'
$ws.Range('E13').Value = '  +10.18%  '

$ws.Range('D14').Value = '0.104'
$ws.Range('E14').Value = '  +1.08%  '

$ws.Range('E15').Value = '  +0.38%  '

$ws.Range('D16').Value = '0.870'
$ws.Range('E16').Value = '  +5.50%  '

$ws.Range('D17').Value = '14.51'
$ws.Range('E17').Value = '  +0.32%  '

$ws.Range('D18').Value = '2.209.81'
$ws.Range('E18').Value = '  +0.26%  '

$ws.Range('D19').Value = '41.197.57'
$ws.Range('E19').Value = '  +1.31%  '

$ws.Range('D20').Value = '0.0₃0952'
$ws.Range('E20').Value = '  +2.31%  '

$ws.Range('D21').Value = '6.18'
$ws.Range('E21').Value = '  +3.12%  '

$ws.Range('D22').Value = '71.92'
$ws.Range('E22').Value = '  +0.13%  '

$ws.Range('D23').Value = '232.27'
$ws.Range('E23').Value = '  +1.37%  '

$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  +2.86%  '

$ws.Range('E25').Value = '  +9.12%  '

$ws.Range('D26').Value = '11.83'
$ws.Range('E26').Value = '  +22.99%  '

$ws.Range('E27').Value = '  +0.00%  '

$ws.Range('D28').Value = '2.52'
$ws.Range('E28').Value = '  +6.12%  '

$ws.Range('E29').Value = '  -0.28%  '

$ws.Range('D30').Value = '169.18'
$ws.Range('E30').Value = '  +0.40%  '

$ws.Range('D31').Value = '20.59'
$ws.Range('E31').Value = '  +2.89%  '

$ws.Range('D32').Value = '0.118'
$ws.Range('E32').Value = '  +1.87%  '

$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').Value = '0.123'
$ws.Range('E33').Value = '  +1.05%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.0741'
$ws.Range('E34').Value = '  +6.94%  '

$ws.Range('D35').Value = '5.46'
$ws.Range('E35').Value = '  +8.13%  '

$ws.Range('D36').Value = '27.75'
$ws.Range('E36').Value = '  +22.47%  '

$ws.Range('D37').Value = '4.20'
$ws.Range('E37').Value = '  +11.67%  '

$ws.Range('D38').Value = '4.61'
$ws.Range('E38').Value = '  +2.02%  '

$ws.Range('D39').Value = '0.0299'
$ws.Range('E39').Value = '  +13.88%  '

$ws.Range('D40').Value = '12.75'
$ws.Range('E40').Value = '  +27.71%  '

$ws.Range('D41').Value = '2.20'
$ws.Range('E41').Value = '  -1.19%  '

$ws.Range('D42').Value = '5.70'
$ws.Range('E42').Value = '  -1.22%  '

$ws.Range('D43').Value = '64.17'
$ws.Range('E43').Value = '  +0.75%  '

$ws.Range('D44').Value = '4.99'
$ws.Range('E44').Value = '  +4.48%  '

$ws.Range('D45').Value = '0.201'
$ws.Range('E45').Value = '  +5.98%  '

$ws.Range('D46').Value = '8.59'
$ws.Range('E46').Value = '  +0.16%  '

$ws.Range('E47').Value = '  +3.92%  '

$ws.Range('D49').Value = '1.13'
$ws.Range('E49').Value = '  +5.64%  '

$ws.Range('E50').Value = '  +1.55%  '

$ws.Range('D51').Value = '4.28'
$ws.Range('E51').Value = '  -2.96%  '

